# Mise à jour de l'application
# Append new training-session rows (J-2, 2025-11-20) to the bottom of the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 885

$data = @(
    @("Entrainement", 45981, "Global", "J-2", "Malik Boussaid",     "right back",      "01:01:21", 3.69, 0.52, 3.17, 0.18, 0.11, 0.15, 0.08, 7, 2.96, 32.77, 5.39, 13, 7, 16, 4),
    @("Entrainement", 45981, "Global", "J-2", "Karim Belmahi",      "left forward",    "01:01:52", 3.52, 0.32, 3.2,  0.16, 0.13, 0.03, 0,    3, 2.9,  27.46, 3.83, 11, 0, 4,  0),
    @("Entrainement", 45981, "Global", "J-2", "Karahali Souaré",    "right forward",   "01:01:46", 3.99, 0.43, 3.56, 0.19, 0.14, 0.09, 0.01, 7, 3.52, 30.95, 4.6,  29, 5, 24, 6),
    @("Entrainement", 45981, "Global", "J-2", "Hedi Nasri",         "right back",      "00:35:54", 2.6,  0.48, 2.11, 0.17, 0.1,  0.11, 0.11, 7, 3.74, 34.27, 4.27, 16, 3, 8,  4),
    @("Entrainement", 45981, "Global", "J-2", "Emmanuel Valey",     "left forward",    "00:52:13", 3.12, 0.11, 3.01, 0.11, 0,    0,    0,    0, 2.7,  20.06, 4.04, 11, 1, 7,  1),
    @("Entrainement", 45981, "Global", "J-2", "Levy Ndoutoume",     "left back",       "01:00:14", 3.17, 0.46, 2.7,  0.19, 0.1,  0.12, 0.06, 7, 2.47, 33.08, 4.31, 22, 2, 10, 5),
    @("Entrainement", 45981, "Global", "J-2", "Mattheo Haon",       "right back",      "01:01:15", 3.65, 0.63, 3.02, 0.31, 0.18, 0.11, 0.02, 6, 3.51, 31.07, 4.31, 12, 3, 9,  0),
    @("Entrainement", 45981, "Global", "J-2", "Ilan Ihaddadene",    "center midfield", "01:01:45", 3.59, 0.49, 3.09, 0.26, 0.16, 0.07, 0,    4, 3.39, 27.58, 4.17, 14, 2, 2,  0),
    @("Entrainement", 45981, "Global", "J-2", "Omar Benyounes",     "center midfield", "00:59:57", 3.47, 0.43, 3.02, 0.28, 0.1,  0.06, 0,    7, 2.97, 30.07, 4.67, 22, 7, 9,  3),
    @("Entrainement", 45981, "Global", "J-2", "Kamal Bafounta",     "center midfield", "00:36:12", 2.28, 0.47, 1.8,  0.17, 0.2,  0.11, 0,    5, 3.64, 29.54, 3.75, 5,  0, 2,  0),
    @("Entrainement", 45981, "Global", "J-2", "Jeremie Laurent",    "left forward",    "00:59:51", 3.1,  0.54, 2.56, 0.23, 0.11, 0.16, 0.05, 7, 3.02, 31.96, 4.32, 18, 2, 2,  0),
    @("Entrainement", 45981, "Global", "J-2", "Fareh Wael",         "center midfield", "00:36:12", 2.38, 0.44, 1.93, 0.16, 0.16, 0.13, 0,    7, 3.83, 28.33, 4.24, 10, 1, 3,  1)
)

$lastRow = $startRow + $data.Count - 1

# Column B (dates) uses a specific date style in the existing sheet (style index 1,
# numFmtId 14). Copy that formatting down onto the new date cells so no new,
# redundant number-format/style entries are introduced.
$ws.Cells.Item($startRow - 1, 2).Copy() | Out-Null
$ws.Range("B$($startRow):B$($lastRow)").PasteSpecial(-4122) | Out-Null

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]          # A - Type
    $ws.Cells.Item($r, 2).Value = $row[1]          # B - Date
    $ws.Cells.Item($r, 3).Value = $row[2]          # C - Période
    $ws.Cells.Item($r, 4).Value = $row[3]          # D - MD (J-2)
    $ws.Cells.Item($r, 5).Value = $row[4]          # E - Nom du joueur
    $ws.Cells.Item($r, 6).Value = $row[5]          # F - Poste
    $ws.Cells.Item($r, 7).Value = $row[6]          # G - Temps joué
    $ws.Cells.Item($r, 8).Value = $row[7]          # H - Distance (km)
    $ws.Cells.Item($r, 9).Value = $row[8]          # I - Distance HID
    $ws.Cells.Item($r, 10).Value = $row[9]         # J - Distance 0-15
    $ws.Cells.Item($r, 11).Value = $row[10]        # K - Distance 15-20
    $ws.Cells.Item($r, 12).Value = $row[11]        # L - Distance 20-25
    $ws.Cells.Item($r, 13).Value = $row[12]        # M - Distance 25-30
    $ws.Cells.Item($r, 14).Value = $row[13]        # N - Distance >30
    $ws.Cells.Item($r, 15).Value = $row[14]        # O - # Sprints
    $ws.Cells.Item($r, 16).Value = $row[15]        # P - Vitesse moy.
    $ws.Cells.Item($r, 17).Value = $row[16]        # Q - Vitesse max
    $ws.Cells.Item($r, 18).Value = $row[17]        # R - Accélération max
    $ws.Cells.Item($r, 19).Value = $row[18]        # S - # Acc >3
    $ws.Cells.Item($r, 20).Value = $row[19]        # T - # Acc >4
    $ws.Cells.Item($r, 21).Value = $row[20]        # U - # Dec >3
    $ws.Cells.Item($r, 22).Value = $row[21]        # V - # Dec >4
    $r = $r + 1
}

$window = $ws.Application.ActiveWindow
$window.ScrollRow = 880
$ws.Range("F895").Select() | Out-Null
